$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the completion date (2023-01-23, serial 44949) for rows 9-15 in
# column C, matching the existing date format already applied to these cells.
$date = (Get-Date -Year 2023 -Month 1 -Day 23).Date
foreach ($r in 9..15) {
    $ws.Range("C$r").Value = $date
}

# Update the active selection to C9, as recorded in the saved view state.
$ws.Range("C9").Select()
